$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal TEXT value into a cell, even when the text
# looks like a number (e.g. "113.81"), without leaving any numeric
# auto-conversion or stray cell-style behind. We build the text via a
# formula that yields a string, then collapse it down to a static value
# with Paste Special > Values, which keeps the cell as a plain text cell
# (no NumberFormat residue).
$xlPasteValues = -4163
function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial($xlPasteValues)
}

Set-TextValue $ws.Range("D2") "43.666.27"
$ws.Range("E2").Value = '  -0.24%  '
Set-TextValue $ws.Range("D3") "2.281.05"
$ws.Range("E3").Value = '  -0.53%  '
$ws.Range("E4").Value = '  -0.04%  '
Set-TextValue $ws.Range("D5") "113.81"
$ws.Range("E5").Value = '  +10.60%  '
Set-TextValue $ws.Range("D6") "267.73"
$ws.Range("E6").Value = '  -1.10%  '
Set-TextValue $ws.Range("D7") "0.625"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("E9").Value = '  +0.86%  '
Set-TextValue $ws.Range("D10") "48.51"
$ws.Range("E10").Value = '  +5.30%  '
Set-TextValue $ws.Range("D11") "0.0937"
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("E12").Value = '  +9.62%  '
Set-TextValue $ws.Range("D13") "0.107"
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("E14").Value = '  +1.41%  '
Set-TextValue $ws.Range("D15") "2.623.96"
$ws.Range("E15").Value = '  -0.27%  '
Set-TextValue $ws.Range("D16") "0.874"
$ws.Range("E16").Value = '  +2.07%  '
Set-TextValue $ws.Range("D17") "2.275.31"
$ws.Range("E17").Value = '  -0.75%  '
Set-TextValue $ws.Range("D18") "43.496.16"
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("E19").Value = '  -1.10%  '
Set-TextValue $ws.Range("D20") "7.08"
$ws.Range("E20").Value = '  +12.77%  '
Set-TextValue $ws.Range("D21") "72.07"
$ws.Range("E21").Value = '  -0.37%  '
Set-TextValue $ws.Range("D22") "2.42"
$ws.Range("E22").Value = '  -3.60%  '
Set-TextValue $ws.Range("D23") "9.94"
$ws.Range("E23").Value = '  +7.04%  '
Set-TextValue $ws.Range("D24") "232.75"
$ws.Range("E24").Value = '  -0.34%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  -0.02%  '
Set-TextValue $ws.Range("D27") "11.57"
$ws.Range("E27").Value = '  +2.58%  '
Set-TextValue $ws.Range("D28") "41.47"
$ws.Range("E28").Value = '  +1.27%  '
$ws.Range("E29").Value = '  -1.62%  '
$ws.Range("E30").Value = '  -1.30%  '
Set-TextValue $ws.Range("D31") "173.59"
$ws.Range("E31").Value = '  -2.39%  '
Set-TextValue $ws.Range("D32") "21.52"
$ws.Range("E32").Value = '  -1.42%  '
$ws.Range("E33").Value = '  +0.52%  '
Set-TextValue $ws.Range("D34") "5.67"
$ws.Range("E34").Value = '  +2.85%  '
$ws.Range("E35").Value = '  +0.30%  '
$ws.Range("E36").Value = '  -4.55%  '
Set-TextValue $ws.Range("D37") "0.0352"
$ws.Range("E37").Value = '  -1.90%  '
$ws.Range("E38").Value = '  -3.70%  '
$ws.Range("E39").Value = '  +5.34%  '
Set-TextValue $ws.Range("D40") "14.65"
$ws.Range("E40").Value = '  +19.71%  '
Set-TextValue $ws.Range("D41") "74.96"
$ws.Range("E41").Value = '  +13.95%  '
$ws.Range("E42").Value = '  +3.95%  '
Set-TextValue $ws.Range("D44") "6.25"
$ws.Range("E44").Value = '  +18.37%  '
$ws.Range("E45").Value = '  +0.04%  '
Set-TextValue $ws.Range("D46") "1.38"
$ws.Range("E46").Value = '  +0.36%  '
$ws.Range("E47").Value = '  -1.33%  '
$ws.Range("E48").Value = '  +2.04%  '
Set-TextValue $ws.Range("D49") "102.09"
$ws.Range("E49").Value = '  +3.00%  '
Set-TextValue $ws.Range("D50") "0.0997"
$ws.Range("E50").Value = '  -2.09%  '
Set-TextValue $ws.Range("D51") "0.454"
$ws.Range("E51").Value = '  +2.76%  '

$excel.CutCopyMode = $false

Write-Output "Done updating cryptos list"
